$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 97, shifting the existing
# rows 97-105 down to 99-107 (values for those rows stay identical).
$ws.Rows("97:98").Insert()

# Populate new row 97 with a "Clementina" record dated 2022-07-06 (44748)
$ws.Range("A97").Value2 = 1
$ws.Range("B97").Value2 = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C97").Value2 = 'Arica y Parinacota'
$ws.Range("D97").Value2 = 44748
$ws.Range("E97").Value2 = 15
$ws.Range("F97").Value2 = 'Fruta'
$ws.Range("G97").Value2 = 100102
$ws.Range("H97").Value2 = 'Cítricos'
$ws.Range("I97").Value2 = 100102004
$ws.Range("J97").Value2 = 'Mandarina'
$ws.Range("K97").Value2 = 'Clementina'
$ws.Range("L97").Value2 = 'Segunda'
$ws.Range("M97").Value2 = 250
$ws.Range("N97").Value2 = 14000
$ws.Range("O97").Value2 = 15000
$ws.Range("P97").Value2 = 14500
$ws.Range("Q97").Value2 = '$/caja 20 kilos'
$ws.Range("R97").Value2 = 'Región de Coquimbo'
$ws.Range("S97").Value2 = 725
$ws.Range("T97").Value2 = 20

# Populate new row 98 with a "Clemenuless" record, also dated 2022-07-06 (44748)
$ws.Range("A98").Value2 = 1
$ws.Range("B98").Value2 = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C98").Value2 = 'Arica y Parinacota'
$ws.Range("D98").Value2 = 44748
$ws.Range("E98").Value2 = 15
$ws.Range("F98").Value2 = 'Fruta'
$ws.Range("G98").Value2 = 100102
$ws.Range("H98").Value2 = 'Cítricos'
$ws.Range("I98").Value2 = 100102004
$ws.Range("J98").Value2 = 'Mandarina'
$ws.Range("K98").Value2 = 'Clemenuless'
$ws.Range("L98").Value2 = 'Segunda'
$ws.Range("M98").Value2 = 300
$ws.Range("N98").Value2 = 14000
$ws.Range("O98").Value2 = 15000
$ws.Range("P98").Value2 = 14500
$ws.Range("Q98").Value2 = '$/caja 20 kilos'
$ws.Range("R98").Value2 = 'Región de Coquimbo'
$ws.Range("S98").Value2 = 725
$ws.Range("T98").Value2 = 20
